$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Leading apostrophe forces Excel to store the value as literal text so
# formats like trailing zeros ('3.380') and percents ('0.52%') survive;
# the Style reset afterward keeps the cell on the default (unstyled) format,
# matching the original workbook formatting.
$updates = [ordered]@{
    "D2" = "329.42"
    "E2" = "0.52%"
    "D3" = "44.22"
    "E3" = "0.47%"
    "D4" = "5.508"
    "E4" = "-0.95%"
    "D5" = "0.08018"
    "E5" = "-0.36%"
    "D6" = "2.062"
    "E6" = "7.90%"
    "D7" = "2.628"
    "E7" = "3.42%"
    "D8" = "0.9546"
    "E8" = "1.05%"
    "D9" = "0.1140"
    "E9" = "-2.19%"
    "D10" = "0.1878"
    "E10" = "1.80%"
    "D11" = "10.32"
    "E11" = "7.57%"
    "D12" = "0.09846"
    "E12" = "1.43%"
    "D13" = "0.04865"
    "E13" = "10.95%"
    "D14" = "0.1062"
    "E14" = "-0.48%"
    "D15" = "0.001257"
    "E15" = "-1.89%"
    "E16" = "-2.97%"
    "D17" = "0.006041"
    "E17" = "1.66%"
    "D18" = "3.380"
    "E18" = "-1.09%"
    "D19" = "4.401"
    "E19" = "2.60%"
    "E20" = "-2.57%"
    "D21" = "0.1383"
    "E21" = "1.55%"
    "D22" = "0.2579"
    "E22" = "-2.73%"
    "D23" = "0.001303"
    "E23" = "4.69%"
    "D24" = "0.004360"
    "E24" = "-2.74%"
    "E25" = "-6.37%"
    "D26" = "0.0003747"
    "E26" = "-6.17%"
    "D38" = "0.02575"
    "E38" = "-2.53%"
    "D39" = "0.05764"
    "E39" = "5.14%"
    "D40" = "0.007595"
    "E40" = "0.09%"
    "D41" = "0.1403"
    "E41" = "0.45%"
    "D42" = "0.007322"
    "E42" = "-9.18%"
    "D43" = "0.001993"
    "E43" = "-0.84%"
    "D44" = "0.009044"
    "E44" = "4.86%"
    "D45" = "0.00007002"
    "E45" = "1.42%"
    "D46" = "0.00000000750"
    "E46" = "-0.07%"
    "D47" = "0.0005805"
    "E47" = "-0.12%"
    "D48" = "0.003500"
    "E48" = "54.05%"
    "E49" = "-31.91%"
    "D50" = "0.00002100"
    "E50" = "-0.07%"
    "D51" = "0.0002000"
    "E51" = "-0.07%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}
